# "Check Downloaded Excel File" test fixture update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Location" value for the sample row changed.
$ws.Range("C2").Value = "Spacetoon"

# The "Age" value for the sample row changed.
$ws.Range("B2").Value = 10

# Resize column C ("Location") to fit its new, shorter content.
$ws.Columns("C:C").AutoFit()
$ws.Columns("C:C").ColumnWidth = 8.88671875

# Leave the cursor parked on C8, like the saved fixture does.
$ws.Range("C8").Select()
